# Trn_0322ConsultarTasas.xlsx - "Inicio consulta de costo V!"
# Row 3 of the "Datos" sheet is switched from an "Acierto" (success) test
# case to an "Error" test case: orientacion -> Error, codigoError -> 030,
# resultadoEsperado -> ERROR. The active selection also moves to J8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Row 3 (third data row): orientacion / codigoError / resultadoEsperado
$ws.Range("G3").Value = "Error"
$ws.Range("H3").Value = "030"
$ws.Range("J3").Value = "ERROR"

# Move / update the selection shown when the workbook is reopened.
$ws.Range("J8").Select()
